$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "1.004", "30.577.50").
# Force text format first so Excel does not auto-convert/round them to numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '30.577.50'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').Value = '2.101.31'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.91%  '
$ws.Range('D5').Value = '330.85'
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').Value = '1.005'
$ws.Range('E6').Value = '  +1.09%  '
$ws.Range('D7').Value = '0.5231'
$ws.Range('E7').Value = '  +1.63%  '
$ws.Range('D8').Value = '0.4385'
$ws.Range('E8').Value = '  +2.13%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '53.82'
$ws.Range('E9').Value = '  +23.77%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.08899'
$ws.Range('E10').Value = '  +1.14%  '
$ws.Range('D11').Value = '1.163'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('D12').Value = '24.61'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('D13').Value = '2.103.26'
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('D14').Value = '6.728'
$ws.Range('E14').Value = '  +0.00%  '
$ws.Range('D15').Value = '7.731'
$ws.Range('E15').Value = '  +2.42%  '
$ws.Range('D16').Value = '96.30'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').Value = '1.005'
$ws.Range('E17').Value = '  +1.02%  '
$ws.Range('D18').Value = '0.00001126'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('D19').Value = '0.06624'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').Value = '19.31'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').Value = '1.004'
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('D22').Value = '6.297'
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').Value = '30.610.33'
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').Value = '12.23'
$ws.Range('E24').Value = '  +2.93%  '
$ws.Range('D25').Value = '2.348'
$ws.Range('E25').Value = '  +3.97%  '
$ws.Range('D26').Value = '2.347.41'
$ws.Range('E26').Value = '  +2.06%  '
$ws.Range('D27').Value = '22.38'
$ws.Range('E27').Value = '  -1.54%  '
$ws.Range('D28').Value = '2.615'
$ws.Range('E28').Value = '  +4.13%  '
$ws.Range('D29').Value = '162.78'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').Value = '132.42'
$ws.Range('E30').Value = '  -0.66%  '
$ws.Range('D31').Value = '1.209'
$ws.Range('E31').Value = '  +4.50%  '
$ws.Range('D32').Value = '0.1070'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('D33').Value = '1.671'
$ws.Range('E33').Value = '  +14.43%  '
$ws.Range('D34').Value = '6.211'
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('D35').Value = '3.930'
$ws.Range('E35').Value = '  +2.74%  '
$ws.Range('D36').Value = '10.17'
$ws.Range('E36').Value = '  +9.13%  '
$ws.Range('D37').Value = '0.02583'
$ws.Range('E37').Value = '  +0.80%  '
$ws.Range('D38').Value = '0.06918'
$ws.Range('E38').Value = '  +3.54%  '
$ws.Range('D39').Value = '5.485'
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('D40').Value = '12.65'
$ws.Range('E40').Value = '  +2.05%  '
$ws.Range('D41').Value = '0.2277'
$ws.Range('E41').Value = '  +1.80%  '
$ws.Range('D42').Value = '0.6942'
$ws.Range('E42').Value = '  +2.77%  '
$ws.Range('D43').Value = '1.269'
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('D44').Value = '1.004'
$ws.Range('E44').Value = '  +1.02%  '
$ws.Range('D45').Value = '0.6442'
$ws.Range('E45').Value = '  +3.09%  '
$ws.Range('D46').Value = '13.93'
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('D47').Value = '2.209'
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('B48').Value = 'WEMIXTOKEN'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '1.264'
$ws.Range('E48').Value = '  +13.40%  '
$ws.Range('B49').Value = 'PancakeSwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D49').Value = '3.632'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').Value = '1.251'
$ws.Range('E50').Value = '  -1.00%  '
$ws.Range('D51').Value = '82.32'
$ws.Range('E51').Value = '  -0.40%  '
